$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.297.06'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '''2.609.41'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''518.91'
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").Value = '''149.23'
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = '''0.569'
$ws.Range("E8").Value = '  -4.67%  '
$ws.Range("D9").Value = '''2.610.01'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '''6.36'
$ws.Range("E10").Value = '  -4.22%  '
$ws.Range("D11").Value = '''0.105'
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '''0.342'
$ws.Range("E12").Value = '  -1.04%  '
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '''3.072.24'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").Value = '''60.325.98'
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").Value = '''21.32'
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").Value = '''0.0000138'
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").Value = '''2.624.72'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D19").Value = '''4.62'
$ws.Range("E19").Value = '  -2.62%  '
$ws.Range("D20").Value = '''343.12'
$ws.Range("E20").Value = '  -4.18%  '
$ws.Range("D21").Value = '''10.40'
$ws.Range("E21").Value = '  -1.67%  '
$ws.Range("D22").Value = '''6.10'
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").Value = '''60.59'
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("D25").Value = '''0.418'
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("D26").Value = '''0.163'
$ws.Range("E26").Value = '  -1.70%  '
$ws.Range("D27").Value = '''0.990'
$ws.Range("E27").Value = '  -0.93%  '
$ws.Range("D28").Value = '''0.0₃0817'
$ws.Range("E28").Value = '  -2.65%  '
$ws.Range("D29").Value = '''7.06'
$ws.Range("E29").Value = '  -3.42%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  +0.11%  '
$ws.Range("D32").Value = '''5.97'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("D33").Value = '''18.91'
$ws.Range("E33").Value = '  -2.67%  '
$ws.Range("D34").Value = '''150.19'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = '''3.96'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -4.51%  '
$ws.Range("E37").Value = '  -3.46%  '
$ws.Range("D38").Value = '''0.869'
$ws.Range("E38").Value = '  +3.48%  '
$ws.Range("D39").Value = '''36.54'
$ws.Range("E39").Value = '  +0.61%  '
$ws.Range("D40").Value = '''1.44'
$ws.Range("E40").Value = '  -2.62%  '
$ws.Range("D41").Value = '''3.63'
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").Value = '''286.28'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("D43").Value = '''0.625'
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").Value = '''0.998'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '''0.0547'
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''19.61'
$ws.Range("E47").Value = '  +0.22%  '
$ws.Range("D48").Value = '''0.0233'
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = '''10.39'
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''4.71'
$ws.Range("E50").Value = '  -4.78%  '
$ws.Range("D51").Value = '''1.954.85'
$ws.Range("E51").Value = '  -1.30%  '
